# Formed the consolidated report
# Fill in the "Absent" (column H) values that were left blank/incorrect
# while consolidating the attendance report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 1
$ws.Range("H16").Value = 0
